$wb = $excel.ActiveWorkbook

# Worksheets involved
$ws1 = $wb.Worksheets.Item("Login")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 ("Login"): move the active cell selection to B20 (no longer the selected tab) ---
$ws1.Range("B20").Select()

# --- Sheet2: put the existing shared string "123" (already used on Login!C1) into C1 ---
$ws1.Range("C1").Copy($ws2.Range("C1"))

# --- Sheet2 becomes the active sheet, with its own active cell L14 ---
$ws2.Activate()
$ws2.Range("L14").Select()

$wb.Save()
